$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 entirely (3rd data row is removed)
$ws.Rows.Item(4).Delete()

# Force text number-format on the columns that must stay text even though
# their contents look numeric/date/percentage-like, so Excel's COM layer
# doesn't silently coerce them into numbers/dates/percentages.
$ws.Range("E2:E3").NumberFormat = "@"
$ws.Range("I2:I3").NumberFormat = "@"
$ws.Range("M2:M3").NumberFormat = "@"
$ws.Range("R2:R3").NumberFormat = "@"
$ws.Range("S2:S3").NumberFormat = "@"

# Row 2 - update existing person's data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "sai"
$ws.Range("C2").Value = "j"
$ws.Range("D2").Value = "j"
$ws.Range("E2").Value = "8"
$ws.Range("F2").Value = "saipjaligama@gmail.com"
$ws.Range("G2").Value = "male"
$ws.Range("H2").Value = 23
$ws.Range("I2").Value = "2023-10-25"
$ws.Range("J2").Value = "5035 s east end S2402`ns2402"
$ws.Range("K2").Value = "ILLINOIS"
$ws.Range("L2").Value = "Preferred Plus Non Tobacco"
$ws.Range("M2").Value = "200000"
$ws.Range("N2").Value = "level"
$ws.Range("O2").Value = "monthly"
$ws.Range("P2").Value = "maximum"
$ws.Range("Q2").Value = "ltc_rider"
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = "2%"
$ws.Range("T2").Value = "Preferred Tobacco"
$ws.Range("U2").Value = "10-Year"

# Row 3 - update existing person's data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Sai"
$ws.Range("C3").Value = "p"
$ws.Range("D3").Value = "p"
$ws.Range("E3").Value = "8155933548"
$ws.Range("F3").Value = "saijaligama@hotmail.com"
$ws.Range("G3").Value = "male"
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = "2023-10-04"
$ws.Range("J3").Value = "13804 Summit Commons BLvd apt I"
$ws.Range("K3").Value = "North Carolina"
$ws.Range("L3").Value = "Preferred Plus Non Tobacco"
$ws.Range("M3").Value = "22222"
$ws.Range("N3").Value = "level"
$ws.Range("O3").Value = "monthly"
$ws.Range("P3").Value = "maximum"
$ws.Range("Q3").Value = "ltc_rider"
$ws.Range("R3").Value = "222"
$ws.Range("S3").Value = "2%"
$ws.Range("T3").Value = "Preferred Tobacco"
$ws.Range("U3").Value = "10-Year"
